# "submit the Gantt chart and Project plan files"
#
# Fill in the remaining Actual Start / Actual Duration / % Complete values
# for the activities in the "Project Planner" Gantt chart (rows 16-28),
# nudge a couple of row heights, and update the saved view (zoom + the
# last-selected cell) to reflect where the author left off.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Project Planner")

# Row -> (Actual Start, Actual Duration, % Complete)
$updates = @(
    @(16, 23, 4, 1),
    @(17, 28, 5, 1),
    @(18, 34, 5, 1),
    @(19, 40, 3, 1),
    @(20, 44, 3, 1),
    @(21, 47, 3, 1),
    @(22, 50, 3, 1),
    @(23, 53, 5, 1),
    @(24, 58, 2, 1),
    @(25, 60, 4, 1),
    @(26, 64, 3, 1),
    @(27, 67, 4, 1),
    @(28, 71, 2, 1)
)

foreach ($u in $updates) {
    $row = $u[0]
    $ws.Range("E$row").Value = $u[1]
    $ws.Range("F$row").Value = $u[2]
    $ws.Range("G$row").Value = $u[3]
}

# Minor row height tweaks made while reviewing the chart.
$ws.Rows.Item(3).RowHeight = 39.9
$ws.Rows.Item(29).RowHeight = 0.9

# Leave the workbook zoomed out a bit and with the selection where the
# author last left it.
$excel.ActiveWindow.Zoom = 85
$ws.Range("F31").Select() | Out-Null
